# Change date of task's end
#
# The "Команды архивирования" (compact/decompact) row's "Когда сделана"
# (when-done) cell is empty; fill it in with "25.04" using the same
# direct (Times New Roman / 28 half-points / no w:lang) run formatting
# already used elsewhere in the table. Word also relocates its hidden
# "_GoBack" bookmark (last-edit marker) to sit right after the newly
# typed text, removing it from its previous location (the "01.05" cell
# of the "Локализация" row).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1) "Команды архивирования" row, "Когда сделана" column: was empty,
#        now gets "25.04" plus the relocated _GoBack bookmark. ---
$targetCell = $t.Cell(12, 4)
$targetRange = $targetCell.Range

$newCellXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00E71CAC" w:rsidRDefault="00E71CAC" w:rsidP="004E6E5B">
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t>25.04</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$targetRange.InsertXML($newCellXml) | Out-Null

# --- 2) "Локализация" row, "Дедлайн" column ("01.05"): drop the
#        _GoBack bookmark that used to live here - it moved above. ---
$oldCell = $t.Cell(14, 3)
$oldRange = $oldCell.Range

$oldCellXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00E71CAC" w:rsidRDefault="0025665B" w:rsidP="004E6E5B">
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="ru-RU"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="ru-RU"/>
</w:rPr>
<w:t>01.05</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$oldRange.InsertXML($oldCellXml) | Out-Null

Write-Output "Updated 'Команды архивирования' row end-date to 25.04 and relocated the _GoBack bookmark."
